$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7393161058425903
$ws.Range("B1").Value = 1.060907959938049
$ws.Range("C1").Value = 2.018189668655396
$ws.Range("D1").Value = 3.457665920257568
$ws.Range("E1").Value = 3.480510234832764
